# Changes to improve performance
# Update the "row-3" individual quota value labels on Sheet1 (columns I and K,
# rows 2-11) to "row-4" labels, and move the sheet's active selection/view to
# column K (D1 as top-left, K3:K11 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 9).Value  = "individualQuotaValuesId-row-4-totalQuota"   # column I
    $ws.Cells.Item($r, 11).Value = "individualQuotaValuesId-row-4-colorQuota"   # column K
}

$ws.Activate()
$ws.Range("K3:K11").Select()
$excel.ActiveWindow.ScrollColumn = 4
